$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells in D2:E51 stay as text (avoid Excel auto-converting
# numeric-looking strings like "1.24" or "37.459.95" into numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '37.459.95'
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = '2.052.83'
$ws.Range("E3").Value = '  +3.89%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '252.49'
$ws.Range("E5").Value = '  +3.39%  '
$ws.Range("E6").Value = '  +2.43%  '
$ws.Range("D7").Value = '66.39'
$ws.Range("E7").Value = '  +17.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +6.89%  '
$ws.Range("D10").Value = '59.97'
$ws.Range("E10").Value = '  +3.41%  '
$ws.Range("D11").Value = '0.0766'
$ws.Range("E11").Value = '  +4.85%  '
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("D14").Value = '14.96'
$ws.Range("E14").Value = '  +4.92%  '
$ws.Range("D15").Value = '2.355.02'
$ws.Range("E15").Value = '  +4.02%  '
$ws.Range("D16").Value = '21.68'
$ws.Range("E16").Value = '  +23.89%  '
$ws.Range("E17").Value = '  +6.44%  '
$ws.Range("D18").Value = '2.052.03'
$ws.Range("E18").Value = '  +3.96%  '
$ws.Range("D19").Value = '37.366.50'
$ws.Range("E19").Value = '  +5.54%  '
$ws.Range("D20").Value = '73.66'
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").Value = '0.0₃0877'
$ws.Range("E21").Value = '  +4.56%  '
$ws.Range("D22").Value = '5.44'
$ws.Range("E22").Value = '  +6.15%  '
$ws.Range("D23").Value = '239.99'
$ws.Range("E23").Value = '  +3.40%  '
$ws.Range("E24").Value = '  +5.15%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +4.79%  '
$ws.Range("D27").Value = '9.81'
$ws.Range("E27").Value = '  +8.44%  '
$ws.Range("D28").Value = '161.61'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").Value = '20.05'
$ws.Range("E30").Value = '  +28.61%  '
$ws.Range("D31").Value = '5.26'
$ws.Range("E31").Value = '  +8.79%  '
$ws.Range("E32").Value = '  +3.46%  '
$ws.Range("E33").Value = '  +8.08%  '
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +10.18%  '
$ws.Range("E35").Value = '  +5.79%  '
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("E37").Value = '  +4.40%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").Value = '6.05'
$ws.Range("E39").Value = '  +18.23%  '
$ws.Range("E40").Value = '  +34.52%  '
$ws.Range("E41").Value = '  +17.23%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '3.01'
$ws.Range("E42").Value = '  +4.63%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.24'
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("E44").Value = '  +6.28%  '
$ws.Range("E45").Value = '  +4.14%  '
$ws.Range("D46").Value = '17.03'
$ws.Range("E46").Value = '  +7.58%  '
$ws.Range("E47").Value = '  +7.05%  '
$ws.Range("D48").Value = '95.80'
$ws.Range("E48").Value = '  +5.35%  '
$ws.Range("D49").Value = '1.422.57'
$ws.Range("E49").Value = '  +3.33%  '
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").Value = '46.75'
$ws.Range("E51").Value = '  +2.00%  '

# Restore original (default) cell formatting/style so only the values differ,
# matching the source diff which contains no style changes.
$ws.Range("D2:E51").ClearFormats()

Write-Output "done"
